{"js": "// Update the division-problem table: each \"NN\u00f7N=\" prompt is replaced by a\n// new prompt, keyed by its (row, column) position in the single table so\n// that duplicate prompt text (e.g. \"91\u00f75=\" appears twice) is handled\n// correctly.\nconst replacements = [\n  { row: 0, col: 0, oldText: \"23\u00f79=\", newText: \"17\u00f78=\" },\n  { row: 0, col: 1, oldText: \"70\u00f79=\", newText: \"94\u00f72=\" },\n  { row: 0, col: 2, oldText: \"30\u00f78=\", newText: \"34\u00f77=\" },\n  { row: 0, col: 3, oldText: \"63\u00f78=\", newText: \"92\u00f76=\" },\n  { row: 0, col: 4, oldText: \"91\u00f75=\", newText: \"67\u00f72=\" },\n  { row: 4, col: 0, oldText: \"65\u00f73=\", newText: \"57\u00f72=\" },\n  { row: 4, col: 1, oldText: \"92\u00f75=\", newText: \"94\u00f78=\" },\n  { row: 4, col: 2, oldText: \"16\u00f78=\", newText: \"23\u00f74=\" },\n  { row: 4, col: 3, oldText: \"66\u00f74=\", newText: \"56\u00f78=\" },\n  { row: 4, col: 4, oldText: \"58\u00f72=\", newText: \"91\u00f73=\" },\n  { row: 8, col: 0, oldText: \"88\u00f79=\", newText: \"45\u00f74=\" },\n  { row: 8, col: 1, oldText: \"91\u00f75=\", newText: \"15\u00f79=\" },\n  { row: 8, col: 2, oldText: \"14\u00f75=\", newText: \"78\u00f79=\" },\n  { row: 8, col: 3, oldText: \"85\u00f76=\", newText: \"10\u00f77=\" },\n  { row: 8, col: 4, oldText: \"71\u00f74=\", newText: \"36\u00f73=\" },\n  { row: 12, col: 0, oldText: \"61\u00f78=\", newText: \"98\u00f78=\" },\n  { row: 12, col: 1, oldText: \"71\u00f76=\", newText: \"39\u00f79=\" },\n  { row: 12, col: 2, oldText: \"50\u00f74=\", newText: \"34\u00f77=\" },\n  { row: 12, col: 3, oldText: \"80\u00f73=\", newText: \"41\u00f76=\" },\n  { row: 12, col: 4, oldText: \"99\u00f75=\", newText: \"85\u00f77=\" },\n  { row: 16, col: 0, oldText: \"99\u00f78=\", newText: \"25\u00f78=\" },\n  { row: 16, col: 1, oldText: \"70\u00f72=\", newText: \"53\u00f76=\" },\n  { row: 16, col: 2, oldText: \"86\u00f77=\", newText: \"88\u00f77=\" },\n  { row: 16, col: 3, oldText: \"82\u00f78=\", newText: \"91\u00f75=\" },\n  { row: 16, col: 4, oldText: \"77\u00f72=\", newText: \"47\u00f78=\" },\n];\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\n\nfor (const r of replacements) {\n  const cell = table.getCell(r.row, r.col);\n  const paragraphs = cell.body.paragraphs;\n  paragraphs.load(\"items\");\n  await context.sync();\n\n  const paragraph = paragraphs.items[0];\n  paragraph.load(\"text\");\n  await context.sync();\n\n  if (paragraph.text !== r.oldText) {\n    throw new Error(\n      `Unexpected text at row ${r.row}, col ${r.col}: expected \"${r.oldText}\", found \"${paragraph.text}\"`\n    );\n  }\n\n  // insertText with Replace keeps the paragraph/run formatting (font,\n  // size, justification) intact, unlike clearing + re-inserting the body.\n  paragraph.insertText(r.newText, Word.InsertLocation.replace);\n}\n\nawait context.sync();\n", "ps1": "# Update the division-problem table: each \"NN\u00f7N=\" prompt is replaced by a\n# new prompt, keyed by its (row, column) position in the single table so\n# that duplicate prompt text (e.g. \"91\u00f75=\" appears twice) is handled\n# correctly.\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n$replacements = @(\n    @{ Row = 1;  Col = 1; Old = \"23\u00f79=\"; New = \"17\u00f78=\" },\n    @{ Row = 1;  Col = 2; Old = \"70\u00f79=\"; New = \"94\u00f72=\" },\n    @{ Row = 1;  Col = 3; Old = \"30\u00f78=\"; New = \"34\u00f77=\" },\n    @{ Row = 1;  Col = 4; Old = \"63\u00f78=\"; New = \"92\u00f76=\" },\n    @{ Row = 1;  Col = 5; Old = \"91\u00f75=\"; New = \"67\u00f72=\" },\n    @{ Row = 5;  Col = 1; Old = \"65\u00f73=\"; New = \"57\u00f72=\" },\n    @{ Row = 5;  Col = 2; Old = \"92\u00f75=\"; New = \"94\u00f78=\" },\n    @{ Row = 5;  Col = 3; Old = \"16\u00f78=\"; New = \"23\u00f74=\" },\n    @{ Row = 5;  Col = 4; Old = \"66\u00f74=\"; New = \"56\u00f78=\" },\n    @{ Row = 5;  Col = 5; Old = \"58\u00f72=\"; New = \"91\u00f73=\" },\n    @{ Row = 9;  Col = 1; Old = \"88\u00f79=\"; New = \"45\u00f74=\" },\n    @{ Row = 9;  Col = 2; Old = \"91\u00f75=\"; New = \"15\u00f79=\" },\n    @{ Row = 9;  Col = 3; Old = \"14\u00f75=\"; New = \"78\u00f79=\" },\n    @{ Row = 9;  Col = 4; Old = \"85\u00f76=\"; New = \"10\u00f77=\" },\n    @{ Row = 9;  Col = 5; Old = \"71\u00f74=\"; New = \"36\u00f73=\" },\n    @{ Row = 13; Col = 1; Old = \"61\u00f78=\"; New = \"98\u00f78=\" },\n    @{ Row = 13; Col = 2; Old = \"71\u00f76=\"; New = \"39\u00f79=\" },\n    @{ Row = 13; Col = 3; Old = \"50\u00f74=\"; New = \"34\u00f77=\" },\n    @{ Row = 13; Col = 4; Old = \"80\u00f73=\"; New = \"41\u00f76=\" },\n    @{ Row = 13; Col = 5; Old = \"99\u00f75=\"; New = \"85\u00f77=\" },\n    @{ Row = 17; Col = 1; Old = \"99\u00f78=\"; New = \"25\u00f78=\" },\n    @{ Row = 17; Col = 2; Old = \"70\u00f72=\"; New = \"53\u00f76=\" },\n    @{ Row = 17; Col = 3; Old = \"86\u00f77=\"; New = \"88\u00f77=\" },\n    @{ Row = 17; Col = 4; Old = \"82\u00f78=\"; New = \"91\u00f75=\" },\n    @{ Row = 17; Col = 5; Old = \"77\u00f72=\"; New = \"47\u00f78=\" }\n)\n\nforeach ($r in $replacements) {\n    $cell = $t.Cell($r.Row, $r.Col)\n    $rng = $cell.Range\n    # Trim the trailing cell-end/paragraph marks so only the visible text\n    # is compared/replaced, keeping the run/paragraph formatting intact.\n    [void]$rng.MoveEnd(1, -1)\n    if ($rng.Text -ne $r.Old) {\n        Write-Output \"Unexpected text at row $($r.Row), col $($r.Col): expected '$($r.Old)', found '$($rng.Text)'\"\n    }\n    $rng.Text = $r.New\n}\n"}
